# Update "想去人数" (number of people interested) figures that changed
# in the regenerated gh-pages data dump (commit 456a3b4).
#
# The same events are listed both on the "展览" sheet and the "全部类型"
# (all types) aggregate sheet, so each value needs to be updated in both
# places, at their respective row positions.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAll        = $wb.Worksheets.Item("全部类型")

# Row -> new F-column value, on the "展览" sheet
$exhibitionUpdates = @{
    2  = 9025
    11 = 4025
    20 = 1458
    22 = 536
    27 = 79
    28 = 1025
    31 = 782
    32 = 81
    39 = 209
    41 = 36
}

foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Row -> new F-column value, on the "全部类型" sheet
$allTypesUpdates = @{
    3  = 9025
    12 = 4025
    26 = 1458
    28 = 536
    34 = 79
    35 = 1025
    37 = 782
    38 = 81
    44 = 209
    46 = 36
}

foreach ($row in $allTypesUpdates.Keys) {
    $sheetAll.Range("F$row").Value = $allTypesUpdates[$row]
}
